# Auto-generated script to update 2023 (column J) violent crime figures
# per the commit "Add data for 2023-10-02"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('J2').Value = 5781
$ws.Range('J3').Value = 6166
$ws.Range('J4').Value = 1337
$ws.Range('J5').Value = 472
$ws.Range('J6').Value = 7876
$ws.Range('J7').Value = 21632

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('J2').Value = 174
$ws.Range('J7').Value = 635
$ws.Range('J8').Value = 1361
$ws.Range('J9').Value = 104
$ws.Range('J10').Value = 153
$ws.Range('J11').Value = 337
$ws.Range('J14').Value = 107
$ws.Range('J16').Value = 84
$ws.Range('J18').Value = 180
$ws.Range('J19').Value = 635
$ws.Range('J23').Value = 205
$ws.Range('J24').Value = 68
$ws.Range('J29').Value = 1204
$ws.Range('J31').Value = 195
$ws.Range('J33').Value = 1003
$ws.Range('J34').Value = 102
$ws.Range('J37').Value = 665
$ws.Range('J41').Value = 138
$ws.Range('J42').Value = 903
$ws.Range('J46').Value = 72
$ws.Range('J48').Value = 256
$ws.Range('J51').Value = 267
$ws.Range('J52').Value = 541
$ws.Range('J54').Value = 420
$ws.Range('J55').Value = 291
$ws.Range('J57').Value = 96
$ws.Range('J60').Value = 128
$ws.Range('J63').Value = 72
$ws.Range('J65').Value = 549
$ws.Range('J66').Value = 65
$ws.Range('J67').Value = 818
$ws.Range('J68').Value = 42
$ws.Range('J72').Value = 89
$ws.Range('J73').Value = 207
$ws.Range('J76').Value = 327
$ws.Range('J77').Value = 163
$ws.Range('J79').Value = 621
$ws.Range('J83').Value = 438
$ws.Range('J84').Value = 183
$ws.Range('J85').Value = 899
$ws.Range('J86').Value = 135
$ws.Range('J89').Value = 288
$ws.Range('J91').Value = 243
$ws.Range('J94').Value = 218
$ws.Range('J95').Value = 319
$ws.Range('J97').Value = 181
$ws.Range('J99').Value = 341
$ws.Range('J101').Value = 21632

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('J2').Value = 40
$ws.Range('J7').Value = 107

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('J3').Value = 190
$ws.Range('J6').Value = 206
$ws.Range('J7').Value = 635

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('J3').Value = 68
$ws.Range('J6').Value = 140
$ws.Range('J7').Value = 337

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('J2').Value = 91
$ws.Range('J3').Value = 79
$ws.Range('J7').Value = 288

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('J2').Value = 237
$ws.Range('J3').Value = 323
$ws.Range('J7').Value = 899

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('J2').Value = 130
$ws.Range('J6').Value = 216
$ws.Range('J7').Value = 541

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('J2').Value = 375
$ws.Range('J3').Value = 415
$ws.Range('J6').Value = 459
$ws.Range('J7').Value = 1361

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('J3').Value = 161
$ws.Range('J7').Value = 438

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('J3').Value = 335
$ws.Range('J5').Value = 42
$ws.Range('J6').Value = 343
$ws.Range('J7').Value = 1003

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('J3').Value = 113
$ws.Range('J7').Value = 319

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('J2').Value = 198
$ws.Range('J3').Value = 222
$ws.Range('J6').Value = 195
$ws.Range('J7').Value = 665

$ws = $wb.Worksheets.Item('New City')
$ws.Range('J2').Value = 160
$ws.Range('J3').Value = 159
$ws.Range('J7').Value = 549

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('J2').Value = 96
$ws.Range('J6').Value = 89
$ws.Range('J7').Value = 341

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('J3').Value = 52
$ws.Range('J7').Value = 195

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('J3').Value = 311
$ws.Range('J6').Value = 217
$ws.Range('J7').Value = 818

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('J6').Value = 55
$ws.Range('J7').Value = 183

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('J6').Value = 203
$ws.Range('J7').Value = 420

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('J3').Value = 421
$ws.Range('J7').Value = 1204

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('J4').Value = 39
$ws.Range('J7').Value = 256

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('J3').Value = 185
$ws.Range('J6').Value = 240
$ws.Range('J7').Value = 635

$ws = $wb.Worksheets.Item('River North')
$ws.Range('J3').Value = 67
$ws.Range('J6').Value = 182
$ws.Range('J7').Value = 327

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('J5').Value = 2
$ws.Range('J6').Value = 75
$ws.Range('J7').Value = 138

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('J2').Value = 195
$ws.Range('J6').Value = 467
$ws.Range('J7').Value = 903

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('J2').Value = 35
$ws.Range('J3').Value = 32
$ws.Range('J7').Value = 153

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('J3').Value = 65
$ws.Range('J6').Value = 145
$ws.Range('J7').Value = 291

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('J4').Value = 11
$ws.Range('J7').Value = 68

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range('J6').Value = 29
$ws.Range('J7').Value = 72

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('J2').Value = 55
$ws.Range('J7').Value = 205

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('J3').Value = 99
$ws.Range('J4').Value = 9
$ws.Range('J7').Value = 243

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('J2').Value = 174
$ws.Range('J3').Value = 213
$ws.Range('J4').Value = 36
$ws.Range('J6').Value = 180
$ws.Range('J7').Value = 621

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range('J6').Value = 85
$ws.Range('J7').Value = 180

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range('J3').Value = 26
$ws.Range('J7').Value = 102

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('J2').Value = 39
$ws.Range('J3').Value = 45
$ws.Range('J7').Value = 218

$ws = $wb.Worksheets.Item('North Center')
$ws.Range('J4').Value = 4
$ws.Range('J7').Value = 65

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range('J6').Value = 34
$ws.Range('J7').Value = 104

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('J2').Value = 71
$ws.Range('J3').Value = 52
$ws.Range('J6').Value = 69
$ws.Range('J7').Value = 207

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('J3').Value = 45
$ws.Range('J7').Value = 174

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('J5').Value = 1
$ws.Range('J6').Value = 124
$ws.Range('J7').Value = 181

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('J6').Value = 25
$ws.Range('J7').Value = 135

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('J2').Value = 64
$ws.Range('J3').Value = 70
$ws.Range('J6').Value = 102
$ws.Range('J7').Value = 267

$ws = $wb.Worksheets.Item('North Park')
$ws.Range('J2').Value = 17
$ws.Range('J7').Value = 42

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range('J3').Value = 25
$ws.Range('J7').Value = 96

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('J3').Value = 37
$ws.Range('J7').Value = 128

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('J6').Value = 32
$ws.Range('J7').Value = 89

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range('J2').Value = 62
$ws.Range('J7').Value = 163

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range('J6').Value = 66
$ws.Range('J7').Value = 84
